# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (column A) uses 4 distinct emoji markers as status
# indicators. Replace them (in place, same row/meaning) with simpler
# text / emoji equivalents:
#   Red book    (📕) -> -3
#   Blue book   (📘) -> warning sign (⚠️)
#   Green book  (📗) -> check mark (✅)
#   Orange book (📙) -> +3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$redBook    = "📕"
$blueBook   = "📘"
$greenBook  = "📗"
$orangeBook = "📙"

$warningSign = "⚠️"
$checkMark   = "✅"

$used = $ws.UsedRange
$rowCount = $used.Row + $used.Rows.Count - 1

# Excel auto-converts number-looking text (e.g. "-3", "+3") to a real
# numeric value when it is written straight into a General-formatted
# cell, which would lose the fact that it must stay a text value.
# Write such values first into a scratch cell that has been explicitly
# formatted as Text, then copy/paste only the *value* from there onto
# the real target cell: this keeps the text stored as text while
# leaving the target cell's own formatting untouched.
$scratch = $ws.Cells.Item($rowCount + 1000, $used.Column + $used.Columns.Count + 10)
$scratch.NumberFormat = "@"

function Set-TextValue($targetCell, $text) {
    $scratch.Value2 = $text
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163)  # xlPasteValues
}

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -eq $redBook) {
        Set-TextValue $cell "-3"
    } elseif ($val -eq $blueBook) {
        $cell.Value2 = $warningSign
    } elseif ($val -eq $greenBook) {
        $cell.Value2 = $checkMark
    } elseif ($val -eq $orangeBook) {
        Set-TextValue $cell "+3"
    }
}

$scratch.Clear()
$excel.CutCopyMode = 0
